# Update "想去人数" (F column) counts across all four sheets to match
# the refreshed data snapshot ("Update gh-pages to output generated at 456a3b4").
# Only the numeric values in column F change; everything else is untouched.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 2436
$ws.Range("F8").Value  = 368
$ws.Range("F9").Value  = 3375
$ws.Range("F10").Value = 892
$ws.Range("F11").Value = 101
$ws.Range("F15").Value = 915
$ws.Range("F18").Value = 410
$ws.Range("F21").Value = 86
$ws.Range("F23").Value = 4039
$ws.Range("F24").Value = 17
$ws.Range("F26").Value = 1182

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F16").Value = 10
$ws.Range("F23").Value = 8
$ws.Range("F24").Value = 9
$ws.Range("F40").Value = 376
$ws.Range("F44").Value = 4
$ws.Range("F48").Value = 310

# Sheet "本地生活" (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value  = 129
$ws.Range("F12").Value = 2865
$ws.Range("F13").Value = 401
$ws.Range("F14").Value = 730
$ws.Range("F15").Value = 77

# Sheet "全部类型" (All Types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value  = 401
$ws.Range("F9").Value  = 730
$ws.Range("F11").Value = 77
$ws.Range("F16").Value = 368
$ws.Range("F17").Value = 3375
$ws.Range("F19").Value = 892
$ws.Range("F20").Value = 101
$ws.Range("F25").Value = 915
$ws.Range("F27").Value = 8
$ws.Range("F28").Value = 9
$ws.Range("F30").Value = 410
$ws.Range("F38").Value = 86
$ws.Range("F41").Value = 4039
$ws.Range("F42").Value = 376
$ws.Range("F48").Value = 1182
